# "Add files via upload" — rename the "Audit Type" column to "Prompt Type"
# and its two values ("IT Audit" -> "IT Prompts", "Business Audit" ->
# "Business Prompts"), then adjust the saved view state (scroll position /
# selection) and give column C an explicit best-fit width, matching the
# refreshed domain-prompts.xlsx upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: rename header + remap the two category values ------------
$ws.Cells.Item(1, 3).Value = "Prompt Type"

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -eq "IT Audit") {
        $cell.Value = "IT Prompts"
    } elseif ($current -eq "Business Audit") {
        $cell.Value = "Business Prompts"
    }
}

# --- Column C width: give it an explicit best-fit-like width ------------
$ws.Columns.Item(3).ColumnWidth = 12

# --- View state: scroll so row 59 is at the top, select D71 -------------
$win = $excel.ActiveWindow
$win.ScrollRow = 59
$win.ScrollColumn = 1
$ws.Range("D71").Select() | Out-Null
